# Add accent6 font color to the review comment about eojeol accuracy / morpheme F1 score
# (item 6) in both the English and Korean text boxes on slide 2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# English content placeholder: "6. The Author should define evaluation parameters "eojeol" accuracy and "morpheme F1" score."
$shEn = $s.Shapes.Item(2)
$trEn = $shEn.TextFrame.TextRange
$paraEn = $trEn.Paragraphs(6)
$rangeEn = $trEn.Characters($paraEn.Start, $paraEn.Length)
$rangeEn.Font.Color.ObjectThemeColor = 10

# Korean content placeholder: "6. 저자는 평가 파라미터인 '어절' 정확도와 '형태소 F1' 점수를 정의해야 합니다."
$shKo = $s.Shapes.Item(3)
$trKo = $shKo.TextFrame.TextRange
$paraKo = $trKo.Paragraphs(6)
$rangeKo = $trKo.Characters($paraKo.Start, $paraKo.Length)
$rangeKo.Font.Color.ObjectThemeColor = 10
